# Apply the cryptos-list price/volume refresh for the GitHub Actions run.
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the source t="inlineStr" cells) instead of
# auto-coercing number-looking strings like "219.15" into floats; the
# Style reset below drops the quote-prefix formatting Excel applies so the
# cell keeps its original (default) style index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}
Set-TextCell "D2" '26.313.32'
Set-TextCell "E2" '  +0.90%  '
Set-TextCell "D3" '1.665.53'
Set-TextCell "E3" '  +0.80%  '
Set-TextCell "E4" '  +0.82%  '
Set-TextCell "D5" '219.15'
Set-TextCell "E5" '  +0.76%  '
Set-TextCell "D6" '0.5356'
Set-TextCell "E6" '  +1.98%  '
Set-TextCell "E7" '  +0.74%  '
Set-TextCell "D8" '0.2664'
Set-TextCell "E8" '  +2.76%  '
Set-TextCell "D9" '0.06404'
Set-TextCell "E9" '  +1.06%  '
Set-TextCell "D10" '20.73'
Set-TextCell "E10" '  +1.92%  '
Set-TextCell "D11" '0.07846'
Set-TextCell "E11" '  +0.58%  '
Set-TextCell "D12" '4.570'
Set-TextCell "E12" '  +1.66%  '
Set-TextCell "D13" '1.668.17'
Set-TextCell "E13" '  +0.80%  '
Set-TextCell "D14" '1.894.21'
Set-TextCell "E14" '  +0.78%  '
Set-TextCell "D15" '0.5536'
Set-TextCell "E15" '  +1.03%  '
Set-TextCell "D16" '0.0₅8228'
Set-TextCell "E16" '  -0.02%  '
Set-TextCell "D17" '65.83'
Set-TextCell "E17" '  +0.70%  '
Set-TextCell "E18" '  +0.76%  '
Set-TextCell "D19" '4.695'
Set-TextCell "E19" '  +2.70%  '
Set-TextCell "D20" '193.95'
Set-TextCell "E20" '  +1.64%  '
Set-TextCell "D21" '10.28'
Set-TextCell "E21" '  +2.27%  '
Set-TextCell "D22" '6.049'
Set-TextCell "E22" '  +0.32%  '
Set-TextCell "E23" '  +0.79%  '
Set-TextCell "D24" '146.26'
Set-TextCell "E24" '  +2.85%  '
Set-TextCell "D25" '0.1232'
Set-TextCell "E25" '  +0.09%  '
Set-TextCell "D26" '7.213'
Set-TextCell "E26" '  -0.31%  '
Set-TextCell "D27" '16.14'
Set-TextCell "E27" '  +0.52%  '
Set-TextCell "E28" '  +4.09%  '
Set-TextCell "D29" '0.05842'
Set-TextCell "E29" '  +0.16%  '
Set-TextCell "D30" '1.283'
Set-TextCell "E30" '  +0.86%  '
Set-TextCell "D31" '3.615'
Set-TextCell "E31" '  +1.98%  '
Set-TextCell "D32" '3.280'
Set-TextCell "E32" '  +0.82%  '
Set-TextCell "E33" '  +2.26%  '
Set-TextCell "D34" '0.9705'
Set-TextCell "E34" '  +2.73%  '
Set-TextCell "D35" '2.824'
Set-TextCell "E36" '  +0.27%  '
Set-TextCell "D37" '0.5826'
Set-TextCell "E37" '  +1.74%  '
Set-TextCell "E38" '  +0.04%  '
Set-TextCell "D39" '0.8753'
Set-TextCell "E39" '  +3.73%  '
Set-TextCell "D40" '5.865'
Set-TextCell "E40" '  +1.94%  '
Set-TextCell "B41" 'Quant'
Set-TextCell "C41" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell "D41" '105.35'
Set-TextCell "E41" '  +1.84%  '
Set-TextCell "B42" 'Maker'
Set-TextCell "C42" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell "D42" '1.052.30'
Set-TextCell "E42" '  +2.64%  '
Set-TextCell "E43" '  +0.74%  '
Set-TextCell "D44" '1.805.05'
Set-TextCell "E44" '  +0.55%  '
Set-TextCell "D45" '57.95'
Set-TextCell "E45" '  +1.66%  '
Set-TextCell "D46" '1.013'
Set-TextCell "E46" '  +1.39%  '
Set-TextCell "B47" 'Mantle'
Set-TextCell "C47" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell "D47" '0.4386'
Set-TextCell "E47" '  +1.64%  '
Set-TextCell "B48" 'BabyDogeCoin'
Set-TextCell "C48" 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell "D48" '0.0₈103'
Set-TextCell "E48" '  -7.67%  '
Set-TextCell "D49" '8.045'
Set-TextCell "E49" '  +3.09%  '
Set-TextCell "D50" '0.05166'
Set-TextCell "E50" '  +0.38%  '
Set-TextCell "D51" '1.415'
Set-TextCell "E51" '  -3.26%  '
